# Update TPM-derived LR-pair metrics (Grn -> Sort1) with recomputed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ "G2"=45.31778466666666; "H2"=135.953354; "I2"=0.1102361023838286; "J2"=0.1102361023838286; "M2"=0.5134303333333333; "N2"=1.540291; "O2"=0.03326489761800302; "P2"=0.03326489761800301; "Q2"=23.26752528733488; "R2"=209.407727586014; "S2"=0.003666992659605758; "T2"=0.003666992659605757 }
    3 = @{ "G3"=45.31778466666666; "H3"=135.953354; "I3"=0.1102361023838286; "J3"=0.1102361023838286; "M3"=1.626140333333333; "N3"=4.878420999999999; "O3"=0.1053568287437347; "P3"=0.1053568287437347; "Q3"=73.69307746378153; "R3"=663.2376971740339; "S3"=0.01161412616022984; "T3"=0.01161412616022984 }
    4 = @{ "G4"=45.31778466666666; "H4"=135.953354; "I4"=0.1102361023838286; "J4"=0.1102361023838286; "M4"=10.254745; "N4"=30.764235; "O4"=0.6643998618255804; "P4"=0.6643998618255803; "Q4"=464.7223257215766; "R4"=4182.50093149419; "S4"=0.07324085119200628; "T4"=0.07324085119200627 }
    5 = @{ "G5"=45.31778466666666; "H5"=135.953354; "I5"=0.1102361023838286; "J5"=0.1102361023838286; "M5"=3.040282666666667; "N5"=9.120848; "O5"=0.1969784118126819; "P5"=0.1969784118126819; "Q5"=137.7788752137991; "R5"=1240.009876924192; "S5"=0.02171413237198676; "T5"=0.02171413237198676 }
    6 = @{ "I6"=0.2429203181515272; "J6"=0.2429203181515272; "M6"=0.5134303333333333; "N6"=1.540291; "O6"=0.03326489761800302; "P6"=0.03326489761800301; "Q6"=51.27317206588077; "R6"=461.4585485929269; "S6"=0.008080719512643273; "T6"=0.008080719512643272 }
    7 = @{ "I7"=0.2429203181515272; "J7"=0.2429203181515272; "M7"=1.626140333333333; "N7"=4.878420999999999; "O7"=0.1053568287437347; "P7"=0.1053568287437347; "S7"=0.025593314357864; "T7"=0.025593314357864 }
    8 = @{ "I8"=0.2429203181515272; "J8"=0.2429203181515272; "M8"=10.254745; "N8"=30.764235; "O8"=0.6643998618255804; "P8"=0.6643998618255803; "Q8"=1024.079160775588; "R8"=9216.712446980295; "S8"=0.1613962258145007; "T8"=0.1613962258145007 }
    9 = @{ "I9"=0.2429203181515272; "J9"=0.2429203181515272; "M9"=3.040282666666667; "N9"=9.120848; "O9"=0.1969784118126819; "P9"=0.1969784118126819; "Q9"=303.6145824982062; "R9"=2732.531242483856; "S9"=0.04785005846651923; "T9"=0.04785005846651923 }
    10 = @{ "G10"=16.49037766666667; "H10"=49.471133; "I10"=0.04011305879538658; "J10"=0.04011305879538658; "M10"=0.5134303333333333; "N10"=1.540291; "O10"=0.03326489761800302; "P10"=0.03326489761800301; "Q10"=8.466660102189222; "R10"=76.199940919703; "S10"=0.00133435679397347; "T10"=0.00133435679397347 }
    11 = @{ "G11"=16.49037766666667; "H11"=49.471133; "I11"=0.04011305879538658; "J11"=0.04011305879538658; "M11"=1.626140333333333; "N11"=4.878420999999999; "O11"=0.1053568287437347; "P11"=0.1053568287437347; "Q11"=26.81566823566589; "R11"=241.341014120993; "S11"=0.004226184665892905; "T11"=0.004226184665892905 }
    12 = @{ "G12"=16.49037766666667; "H12"=49.471133; "I12"=0.04011305879538658; "J12"=0.04011305879538658; "M12"=10.254745; "N12"=30.764235; "O12"=0.6643998618255804; "P12"=0.6643998618255803; "Q12"=169.1046179253617; "R12"=1521.941561328255; "S12"=0.02665111072105622; "T12"=0.02665111072105622 }
    13 = @{ "G13"=16.49037766666667; "H13"=49.471133; "I13"=0.04011305879538658; "J13"=0.04011305879538658; "M13"=3.040282666666667; "N13"=9.120848; "O13"=0.1969784118126819; "P13"=0.1969784118126819; "Q13"=50.13540938675378; "R13"=451.218684480784; "S13"=0.007901406614463978; "T13"=0.007901406614463978 }
    14 = @{ "G14"=249.4253923333333; "H14"=748.276177; "I14"=0.6067305206692575; "J14"=0.6067305206692575; "M14"=0.5134303333333333; "N14"=1.540291; "O14"=0.03326489761800302; "P14"=0.03326489761800301; "Q14"=128.0625623275008; "R14"=1152.563060947507; "S14"=0.02018282865178051; "T14"=0.02018282865178051 }
    15 = @{ "G15"=249.4253923333333; "H15"=748.276177; "I15"=0.6067305206692575; "J15"=0.6067305206692575; "M15"=1.626140333333333; "N15"=4.878420999999999; "O15"=0.1053568287437347; "P15"=0.1053568287437347; "Q15"=405.600690630724; "R15"=3650.406215676516; "S15"=0.06392320355974795; "T15"=0.06392320355974795 }
    16 = @{ "G16"=249.4253923333333; "H16"=748.276177; "I16"=0.6067305206692575; "J16"=0.6067305206692575; "M16"=10.254745; "N16"=30.764235; "O16"=0.6643998618255804; "P16"=0.6643998618255803; "Q16"=2557.793794903288; "R16"=23020.14415412959; "S16"=0.4031116740980171; "T16"=0.4031116740980171 }
    17 = @{ "G17"=249.4253923333333; "H17"=748.276177; "I17"=0.6067305206692575; "J17"=0.6067305206692575; "M17"=3.040282666666667; "N17"=9.120848; "O17"=0.1969784118126819; "P17"=0.1969784118126819; "Q17"=758.3236969375662; "R17"=6824.913272438096; "S17"=0.1195128143597119; "T17"=0.1195128143597119 }
}

foreach ($row in $updates.Keys) {
    foreach ($cellRef in $updates[$row].Keys) {
        $ws.Range($cellRef).Value = $updates[$row][$cellRef]
    }
}